$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.713.55'
$ws.Range('E2').Value = '  +1.27%  '
$ws.Range('D3').Value = '2.489.70'
$ws.Range('E3').Value = '  +1.63%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '532.89'
$ws.Range('E5').Value = '  +3.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.31'
$ws.Range('E6').Value = '  +2.34%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.567'
$ws.Range('E8').Value = '  +2.76%  '
$ws.Range('D9').Value = '2.498.84'
$ws.Range('E9').Value = '  +1.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0989'
$ws.Range('E10').Value = '  +2.83%  '
$ws.Range('E11').Value = '  -2.17%  '
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').Value = '2.934.14'
$ws.Range('E14').Value = '  +1.60%  '
$ws.Range('D15').Value = '58.647.64'
$ws.Range('E15').Value = '  +1.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.21'
$ws.Range('E16').Value = '  +0.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000136'
$ws.Range('E17').Value = '  +1.67%  '
$ws.Range('D18').Value = '2.495.09'
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.57'
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('E20').Value = '  +2.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '319.97'
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.19'
$ws.Range('E22').Value = '  +3.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.92'
$ws.Range('E24').Value = '  +3.93%  '
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('E26').Value = '  +2.46%  '
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.46'
$ws.Range('E28').Value = '  +2.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '172.20'
$ws.Range('E29').Value = '  +2.50%  '
$ws.Range('D30').Value = '0.0₃0754'
$ws.Range('E30').Value = '  +3.01%  '
$ws.Range('E31').Value = '  +4.25%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.18'
$ws.Range('E32').Value = '  +1.34%  '
$ws.Range('B33').Value = 'Aptos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.27'
$ws.Range('E33').Value = '  +1.22%  '
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.10'
$ws.Range('E36').Value = '  +1.54%  '
$ws.Range('E37').Value = '  -3.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.94'
$ws.Range('E38').Value = '  +0.49%  '
$ws.Range('E39').Value = '  +2.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.31'
$ws.Range('E40').Value = '  -0.72%  '
$ws.Range('E41').Value = '  +5.89%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.13'
$ws.Range('E42').Value = '  +2.50%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.46'
$ws.Range('E43').Value = '  +2.20%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '274.31'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '131.05'
$ws.Range('E45').Value = '  +9.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.590'
$ws.Range('E46').Value = '  +0.32%  '
$ws.Range('E47').Value = '  +1.88%  '
$ws.Range('E48').Value = '  +3.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0218'
$ws.Range('E49').Value = '  +3.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.54'
$ws.Range('E50').Value = '  +1.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.74'
$ws.Range('E51').Value = '  +0.22%  '
